# Fruta / hortaliza, semanal
# Insert a new weekly price-report row before row 24 (pushing the existing
# rows 24-40 down to 25-41) and populate the new row with the latest data
# point for "Región de O'Higgins".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(24).Insert()

$ws.Cells.Item(24, 1).Value  = 5
$ws.Cells.Item(24, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(24, 3).Value  = "Maule"
$ws.Cells.Item(24, 4).Value  = 44484
$ws.Cells.Item(24, 5).Value  = 7
$ws.Cells.Item(24, 6).Value  = 100112022
$ws.Cells.Item(24, 7).Value  = "Arveja Verde"
$ws.Cells.Item(24, 8).Value  = "Sin especificar"
$ws.Cells.Item(24, 9).Value  = "Primera"
$ws.Cells.Item(24, 10).Value = 100
$ws.Cells.Item(24, 11).Value = 22000
$ws.Cells.Item(24, 12).Value = 22000
$ws.Cells.Item(24, 13).Value = 22000
$ws.Cells.Item(24, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(24, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(24, 16).Value = 880
$ws.Cells.Item(24, 17).Value = 25
$ws.Cells.Item(24, 18).Value = "Hortaliza"
